# Applies the "DesafiosRPA" configuration commit to Config.xlsx:
#  - Constants!B2: MaxRetryNumber 0 -> 2
#  - Settings!A6:C6: new DesafiosRPAURL row (with hyperlink on B6)
#  - Settings!A8:C8: new DesafiosRPACredential row
#  - Re-selects cells / re-activates the "Settings" sheet so the saved
#    sheet views match the author's final on-screen state.

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets    = $wb.Worksheets.Item("Assets")

# --- Constants sheet: bump MaxRetryNumber from 0 to 2 ---------------------
$constants.Range("B2").Value = 2

# --- Settings sheet: new DesafiosRPA rows ----------------------------------
$settings.Range("A6").Value = "DesafiosRPAURL"
$settings.Range("B6").Value = "https://desafiosrpa.com.br/login.html"
[void]$settings.Hyperlinks.Add($settings.Range("B6"), "https://desafiosrpa.com.br/login.html")
$settings.Range("C6").Value = "Endereço URL do Sistema"

$settings.Range("A8").Value = "DesafiosRPACredential"
$settings.Range("B8").Value = "DesafiosRPA"
$settings.Range("C8").Value = "Nome do Asset no Orchestrator"

# --- Restore the on-screen selections seen in the final workbook ----------
# Select Constants!B3 first (leaves that sheet's cached selection at B3,
# without leaving it the active tab).
[void]$constants.Activate()
[void]$constants.Range("B3").Select()

# Assets keeps its own selection (A2:C2 / active C2) - left untouched so the
# stored selection stays exactly as it was (A2:C2).

# Settings ends up the active/selected tab, with A6 selected.
[void]$settings.Activate()
[void]$settings.Range("A6").Select()
